$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# Fill previously-blank cells (rows 2-13) with the literal text "nan"
# placeholder, mirroring how the source export (pandas/openpyxl) marks
# missing values on every other Card sheet in this workbook.
foreach ($addr in @("D2","E2","F2","G2","H2","I2","J2","K2","L2","M2","N2","O2")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("G3","H3","I3","J3","K3","M3","N3","O3")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","N4","O4")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("D5","H5","J5","K5","M5","N5","O5")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("E6","F6","G6","I6","J6","K6","M6","N6","O6")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("E7","G7","H7","I7","J7","M7","N7","O7")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("E8","H8","I8","J8","K8","M8","N8","O8")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("F9","G9","H9","I9","J9","K9","M9","N9")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("D10","E10","F10","G10","H10","I10","J10","K10","L10","M10","N10","O10")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("D11","E11","F11","G11","H11","I11","J11","K11","L11","M11","N11","O11")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("D12","E12","F12","G12","H12","I12","J12","K12","L12","M12","N12","O12")) { $ws.Range($addr).Value = "nan" }
foreach ($addr in @("B13","C13","D13","E13","F13","G13","H13","I13","J13","K13")) { $ws.Range($addr).Value = "nan" }

# Append the new service event as row 14. Columns B:K have no reading for
# this event (mirrors every other blank-measurement row above), so carry the
# row13 blank formatting down instead of leaving row 14 completely untouched
# - this keeps those cells present-but-empty, same as the rest of the table.
$ws.Range("B13:K13").Copy()
$ws.Range("B14:K14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "9"
$ws.Range("L14").Value = "21\7\2025"
$ws.Range("M14").Value = "159 t"
$ws.Range("N14").Value = "تم سن الفلاتس+ تغيير اول جريده"
$ws.Range("O14").Value = "الخبير"
